# Hindalco prices sheet update (2025-08-12 09:02:43 UTC)
# - Drop the now-superseded first data row (12.08.2025 / 265.25) and let
#   every remaining row shift up by one.
# - Renumber the "Sl.no." column (A) from 21 down to 1 for the 21 rows
#   that remain.
# - The PDF hyperlink that used to live on F2 is gone (it rode along with
#   the deleted row), so clear any leftover hyperlink metadata.
# - Column F was sized to fit the long PDF URL; now that it is unused,
#   shrink it back down to a normal width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete first data row; Excel shifts rows 3..23 up to 2..22
# and updates the used-range dimension automatically.
$ws.Rows.Item(2).Delete()

# Drop all hyperlinks left on the sheet (the one remaining reference is
# stale now that its row moved/merged away).
$ws.Hyperlinks.Delete()

# Renumber column A (Sl.no.) 21 -> 1 across the 21 remaining data rows.
for ($i = 0; $i -lt 21; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = 21 - $i
}

# Shrink column F back from its old 79-char width to 15.
# (ColumnWidth uses Excel's character-width units, which include a fixed
# ~5px/MDW padding offset on top of the stored OOXML "width"; subtracting
# 5/6 compensates for that so the saved width lands exactly on 15.)
$ws.Columns.Item(6).ColumnWidth = 15 - (5 / 6)
